# Scheduled-runner market-data refresh for Fenrir_Profits workbook.
# Updates currentAveragePrice(NQ/HQ), LevePrice(NQ/HQ) and LeveProfit(NQ/HQ)
# columns (H,I,J,K,L,M,N) for leves whose backing item price moved.
# Where an NQ/HQ average price is 0 (no listings), the matching profit
# cell (M for NQ, N for HQ) is cleared rather than written, matching the
# source generator's behaviour for unavailable market data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 74: Adhesive of Antipathy / Wing Glue
$ws.Range("H74").Value = 7214
$ws.Range("I74").Value = 7962.4
$ws.Range("K74").Value = 7962.4
$ws.Range("M74").Value = -7026.4

# Row 77: It's Gonna Grow Back (L) / Wing Glue
$ws.Range("H77").Value = 7214
$ws.Range("I77").Value = 7962.4
$ws.Range("K77").Value = 39812
$ws.Range("M77").Value = -35132

# Row 129: Practical Command / Commanding Craftsman's Draught
$ws.Range("H129").Value = 977.1111
$ws.Range("I129").Value = 743.3333
$ws.Range("K129").Value = 2229.9999
$ws.Range("M129").Value = 2770.0001

$ws = $wb.Worksheets.Item("ARM")
# Row 61: Dealing with the Tough Stuff / Cobalt Ingot
$ws.Range("H61").Value = 3536.7368
$ws.Range("I61").Value = 3757.5151
$ws.Range("K61").Value = 3757.5151
$ws.Range("M61").Value = -3545.5151

# Row 63: Rivets Run through It / Mythrite Rivets
$ws.Range("H63").Value = 1669248.6
$ws.Range("I63").Value = 3032062.8
$ws.Range("J63").Value = 3586.7778
$ws.Range("K63").Value = 3032062.8
$ws.Range("L63").Value = 3586.7778
$ws.Range("M63").Value = -3031376.8
$ws.Range("N63").Value = -4958.7778

# Row 66: A Riveting Revival (L) / Mythrite Rivets
$ws.Range("H66").Value = 1669248.6
$ws.Range("I66").Value = 3032062.8
$ws.Range("J66").Value = 3586.7778
$ws.Range("K66").Value = 15160314
$ws.Range("L66").Value = 17933.889
$ws.Range("M66").Value = -15156882
$ws.Range("N66").Value = -24797.889

# Row 80: A Squire to Inspire / Titanium Hoplon
$ws.Range("H80").Value = 17133.545
$ws.Range("J80").Value = 17133.545
$ws.Range("L80").Value = 17133.545
$ws.Range("N80").Value = -19129.545

# Row 83: All's Fair in Highborn Assassination (L) / Titanium Hoplon
$ws.Range("H83").Value = 17133.545
$ws.Range("J83").Value = 17133.545
$ws.Range("L83").Value = 51400.63499999999
$ws.Range("N83").Value = -61384.63499999999

# Row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Range("H132").Value = 2050351.4
$ws.Range("I132").Value = 2315670.8
$ws.Range("J132").Value = 3601.7144
$ws.Range("K132").Value = 6947012.399999999
$ws.Range("L132").Value = 10805.1432
$ws.Range("M132").Value = -6944482.399999999
$ws.Range("N132").Value = -15865.1432

# Row 136: Metal with Mettle / Cobalt Tungsten Ingot
$ws.Range("H136").Value = 3536.7368
$ws.Range("I136").Value = 3757.5151
$ws.Range("K136").Value = 11272.5453
$ws.Range("M136").Value = -8722.5453

$ws = $wb.Worksheets.Item("BSM")
# Row 82: Spirituality Inspector / Titanium Lump Hammer
$ws.Range("H82").Value = 19638.908
$ws.Range("I82").Value = 11333.333
$ws.Range("J82").Value = 22753.5
$ws.Range("K82").Value = 11333.333
$ws.Range("L82").Value = 22753.5
$ws.Range("M82").Value = -10950.333
$ws.Range("N82").Value = -23519.5

# Row 85: The Clamor for Hammers (L) / Titanium Lump Hammer
$ws.Range("H85").Value = 19638.908
$ws.Range("I85").Value = 11333.333
$ws.Range("J85").Value = 22753.5
$ws.Range("K85").Value = 11333.333
$ws.Range("L85").Value = 22753.5
$ws.Range("M85").Value = -10007.333
$ws.Range("N85").Value = -25405.5

# Row 86: Through Thick and Thin / Adamantite Nugget
$ws.Range("H86").Value = 1406.3334
$ws.Range("I86").Value = 1369.4286
$ws.Range("K86").Value = 1369.4286
$ws.Range("M86").Value = -246.4286

# Row 89: Piercing Eyes Deserve Piercing Shafts (L) / Adamantite Nugget
$ws.Range("H89").Value = 1406.3334
$ws.Range("I89").Value = 1369.4286
$ws.Range("K89").Value = 6847.143
$ws.Range("M89").Value = -1231.143

# Row 134: Ruthenium Supremium / Ruthenium Ingot
$ws.Range("H134").Value = 12207.841
$ws.Range("I134").Value = 12207.841
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 36623.523
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -34088.523
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# Row 62: Splinter in the Sewers / Cedar Lumber
$ws.Range("H62").Value = 5000
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 5000
$ws.Range("N62").Value = -6248
$ws.Range("M62").ClearContents()

# Row 65: The Lumber of Their Discontent (L) / Cedar Lumber
$ws.Range("H65").Value = 5000
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 25000
$ws.Range("N65").Value = -31240
$ws.Range("M65").ClearContents()

# Row 132: Hull Lotta Damage / Ginseng Lumber
$ws.Range("H132").Value = 6413636
$ws.Range("I132").Value = 6945719.5
$ws.Range("J132").Value = 28628
$ws.Range("K132").Value = 20837158.5
$ws.Range("L132").Value = 85884
$ws.Range("M132").Value = -20834628.5
$ws.Range("N132").Value = -90944

$ws = $wb.Worksheets.Item("CUL")
# Row 5: What a Sap / Maple Syrup
$ws.Range("H5").Value = 761.4
$ws.Range("I5").Value = 727.0625
$ws.Range("J5").Value = 898.75
$ws.Range("K5").Value = 2181.1875
$ws.Range("L5").Value = 2696.25
$ws.Range("M5").Value = -2069.1875
$ws.Range("N5").Value = -2920.25

# Row 122: Salt of the North / Northern Sea Salt
$ws.Range("H122").Value = 709.0769
$ws.Range("I122").Value = 530.6842
$ws.Range("J122").Value = 1193.2858
$ws.Range("K122").Value = 4776.1578
$ws.Range("L122").Value = 10739.5722
$ws.Range("M122").Value = -2326.1578
$ws.Range("N122").Value = -15639.5722

# Row 132: More Mezcal / Cooking Mezcal
$ws.Range("H132").Value = 11671.363
$ws.Range("I132").Value = 585.5
$ws.Range("J132").Value = 41233.668
$ws.Range("K132").Value = 5269.5
$ws.Range("L132").Value = 371103.012
$ws.Range("M132").Value = -2739.5
$ws.Range("N132").Value = -376163.012

# Row 135: Not-so-secret Ingredient / Royal Maple Syrup
$ws.Range("H135").Value = 761.4
$ws.Range("I135").Value = 727.0625
$ws.Range("J135").Value = 898.75
$ws.Range("K135").Value = 6543.5625
$ws.Range("L135").Value = 8088.75
$ws.Range("M135").Value = -4008.5625
$ws.Range("N135").Value = -13158.75

$ws = $wb.Worksheets.Item("GSM")
# Row 132: On Board for Lar / Lar Ingot
$ws.Range("H132").Value = 34486680
$ws.Range("I132").Value = 40001210
$ws.Range("J132").Value = 20890.5
$ws.Range("K132").Value = 120003630
$ws.Range("L132").Value = 62671.5
$ws.Range("M132").Value = -120001100
$ws.Range("N132").Value = -67731.5

$ws = $wb.Worksheets.Item("LTW")
# Row 82: Trainin' the Neck / Dragon Leather
$ws.Range("H82").Value = 3760.6843
$ws.Range("I82").Value = 1804.0769
$ws.Range("J82").Value = 8000
$ws.Range("K82").Value = 1804.0769
$ws.Range("L82").Value = 8000
$ws.Range("M82").Value = -1443.0769
$ws.Range("N82").Value = -8722

# Row 85: Training Is Only Skintight (L) / Dragon Leather
$ws.Range("H85").Value = 3760.6843
$ws.Range("I85").Value = 1804.0769
$ws.Range("J85").Value = 8000
$ws.Range("K85").Value = 1804.0769
$ws.Range("L85").Value = 8000
$ws.Range("M85").Value = -556.0769
$ws.Range("N85").Value = -10496

# Row 132: Tenets of Tanning / Silver Lobo Leather
$ws.Range("H132").Value = 4446176.5
$ws.Range("I132").Value = 6897297
$ws.Range("J132").Value = 3520.6875
$ws.Range("K132").Value = 20691891
$ws.Range("L132").Value = 10562.0625
$ws.Range("M132").Value = -20689361
$ws.Range("N132").Value = -15622.0625

$ws = $wb.Worksheets.Item("WVR")
# Row 132: Comfy Cabins / Snow Cotton Cloth
$ws.Range("H132").Value = 4515.3726
$ws.Range("I132").Value = 1216.1136
$ws.Range("J132").Value = 25253.572
$ws.Range("K132").Value = 3648.3408
$ws.Range("L132").Value = 75760.716
$ws.Range("M132").Value = -1118.3408
$ws.Range("N132").Value = -80820.716
